# This script applies a "new weekly data" update to the Betarraga price
# history sheet: a brand-new pair of rows (Primera/Segunda) is inserted
# at the top of the date-ordered data block (rows 274-415), every
# existing pair shifts down by one pair (2 rows), and the oldest pair
# that falls off the bottom is appended as new rows 416-417.
#
# Only columns D (Fecha), J (Volumen), K (Precio minimo), L (Precio
# maximo), M (Precio promedio ponderado) and P (Precio $/Kg) vary from
# row to row in this block; all the other columns (A,B,C,E,F,G,H,I,N,O,Q,R)
# are constant for a given row parity (Primera on even rows, Segunda on
# odd rows) so they do not need to move.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$firstRow = 274
$lastRow = 415

# 1. Snapshot the current D/J/K/L/M/P values for the whole block before
#    we start overwriting anything.
$D = @{}
$J = @{}
$K = @{}
$L = @{}
$M = @{}
$P = @{}

for ($r = $firstRow; $r -le $lastRow; $r++) {
    $D[$r] = $ws.Cells.Item($r, 4).Value2
    $J[$r] = $ws.Cells.Item($r, 10).Value2
    $K[$r] = $ws.Cells.Item($r, 11).Value2
    $L[$r] = $ws.Cells.Item($r, 12).Value2
    $M[$r] = $ws.Cells.Item($r, 13).Value2
    $P[$r] = $ws.Cells.Item($r, 16).Value2
}

# 2. Shift every existing pair down by one pair (2 rows): new row r gets
#    the values that used to live in row r-2. Walk from the bottom up so
#    we never overwrite a source row before it has been read (values are
#    already all cached above, but this keeps the logic self-consistent).
for ($r = $lastRow; $r -ge ($firstRow + 2); $r--) {
    $src = $r - 2
    $ws.Cells.Item($r, 4).Value = $D[$src]
    $ws.Cells.Item($r, 10).Value = $J[$src]
    $ws.Cells.Item($r, 11).Value = $K[$src]
    $ws.Cells.Item($r, 12).Value = $L[$src]
    $ws.Cells.Item($r, 13).Value = $M[$src]
    $ws.Cells.Item($r, 16).Value = $P[$src]
}

# 3. Write the brand-new top pair (newest week of data).
$ws.Cells.Item(274, 4).Value = 44992
$ws.Cells.Item(274, 10).Value = 1700
$ws.Cells.Item(274, 11).Value = 500
$ws.Cells.Item(274, 12).Value = 600
$ws.Cells.Item(274, 13).Value = 550
$ws.Cells.Item(274, 16).Value = 183

$ws.Cells.Item(275, 4).Value = 44992
$ws.Cells.Item(275, 10).Value = 1000
$ws.Cells.Item(275, 11).Value = 400
$ws.Cells.Item(275, 12).Value = 450
$ws.Cells.Item(275, 13).Value = 425
$ws.Cells.Item(275, 16).Value = 142

# 4. Append the pair that fell off the bottom of the block (the original
#    values of rows 414/415) as new rows 416/417, copying the constant
#    columns from the row directly above (414/415) and restoring the
#    date-time number format on column D.
$newRows = @(416, 417)
$sourceRows = @(414, 415)

for ($i = 0; $i -lt 2; $i++) {
    $newRow = $newRows[$i]
    $srcRow = $sourceRows[$i]

    $ws.Cells.Item($newRow, 1).Value = $ws.Cells.Item($srcRow, 1).Value2
    $ws.Cells.Item($newRow, 2).Value = $ws.Cells.Item($srcRow, 2).Value2
    $ws.Cells.Item($newRow, 3).Value = $ws.Cells.Item($srcRow, 3).Value2
    $ws.Cells.Item($newRow, 4).Value = $D[$srcRow]
    $ws.Cells.Item($newRow, 4).NumberFormat = $ws.Cells.Item($srcRow, 4).NumberFormat
    $ws.Cells.Item($newRow, 5).Value = $ws.Cells.Item($srcRow, 5).Value2
    $ws.Cells.Item($newRow, 6).Value = $ws.Cells.Item($srcRow, 6).Value2
    $ws.Cells.Item($newRow, 7).Value = $ws.Cells.Item($srcRow, 7).Value2
    $ws.Cells.Item($newRow, 8).Value = $ws.Cells.Item($srcRow, 8).Value2
    $ws.Cells.Item($newRow, 9).Value = $ws.Cells.Item($srcRow, 9).Value2
    $ws.Cells.Item($newRow, 10).Value = $J[$srcRow]
    $ws.Cells.Item($newRow, 11).Value = $K[$srcRow]
    $ws.Cells.Item($newRow, 12).Value = $L[$srcRow]
    $ws.Cells.Item($newRow, 13).Value = $M[$srcRow]
    $ws.Cells.Item($newRow, 14).Value = $ws.Cells.Item($srcRow, 14).Value2
    $ws.Cells.Item($newRow, 15).Value = $ws.Cells.Item($srcRow, 15).Value2
    $ws.Cells.Item($newRow, 16).Value = $P[$srcRow]
    $ws.Cells.Item($newRow, 17).Value = $ws.Cells.Item($srcRow, 17).Value2
    $ws.Cells.Item($newRow, 18).Value = $ws.Cells.Item($srcRow, 18).Value2
}
